$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2020.0769
$ws.Range("I8").Value = 37.285713
$ws.Range("K8").Value = 111.857139
$ws.Range("M8").Value = 27.142861
$ws.Range("H12").Value = 7930
$ws.Range("I12").Value = 216.66667
$ws.Range("K12").Value = 216.66667
$ws.Range("M12").Value = -46.66667000000001
$ws.Range("H53").Value = 185.15384
$ws.Range("I53").Value = 159.7
$ws.Range("J53").Value = 270
$ws.Range("K53").Value = 159.7
$ws.Range("L53").Value = 270
$ws.Range("M53").Value = 477.3
$ws.Range("N53").Value = -1544
$ws.Range("H70").Value = 3204.1667
$ws.Range("I70").Value = 1787.75
$ws.Range("K70").Value = 5363.25
$ws.Range("M70").Value = -5093.25
$ws.Range("H73").Value = 3204.1667
$ws.Range("I73").Value = 1787.75
$ws.Range("K73").Value = 5363.25
$ws.Range("M73").Value = -4427.25
$ws.Range("H106").Value = 4705.316
$ws.Range("I106").Value = 3855.6667
$ws.Range("J106").Value = 19999
$ws.Range("K106").Value = 3855.6667
$ws.Range("L106").Value = 19999
$ws.Range("M106").Value = -3224.6667
$ws.Range("N106").Value = -21261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 49800.4
$ws.Range("J134").Value = 49800.4
$ws.Range("L134").Value = 49800.4
$ws.Range("N134").Value = -59940.4
$ws.Range("H137").Value = 73647.5
$ws.Range("J137").Value = 73647.5
$ws.Range("L137").Value = 73647.5
$ws.Range("N137").Value = -83847.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 82499
$ws.Range("J60").Value = 82499
$ws.Range("L60").Value = 82499
$ws.Range("N60").Value = -83697
$ws.Range("H86").Value = 2471.7058
$ws.Range("I86").Value = 2902
$ws.Range("K86").Value = 2902
$ws.Range("M86").Value = -1779
$ws.Range("H89").Value = 2471.7058
$ws.Range("I89").Value = 2902
$ws.Range("K89").Value = 14510
$ws.Range("M89").Value = -8894
$ws.Range("H105").Value = 5673.2666
$ws.Range("I105").Value = 3172.7222
$ws.Range("J105").Value = 9424.083000000001
$ws.Range("K105").Value = 3172.7222
$ws.Range("L105").Value = 9424.083000000001
$ws.Range("M105").Value = -1425.7222
$ws.Range("N105").Value = -12918.083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1936.4445
$ws.Range("I31").Value = 1791.125
$ws.Range("K31").Value = 1791.125
$ws.Range("M31").Value = -1496.125
$ws.Range("H34").Value = 1936.4445
$ws.Range("I34").Value = 1791.125
$ws.Range("K34").Value = 1791.125
$ws.Range("M34").Value = -1589.125
$ws.Range("H68").Value = 41998.43
$ws.Range("J68").Value = 41998.43
$ws.Range("L68").Value = 41998.43
$ws.Range("N68").Value = -43496.43
$ws.Range("H71").Value = 41998.43
$ws.Range("J71").Value = 41998.43
$ws.Range("L71").Value = 125995.29
$ws.Range("N71").Value = -133483.29
$ws.Range("H74").Value = 37999
$ws.Range("J74").Value = 37999
$ws.Range("L74").Value = 37999
$ws.Range("N74").Value = -39747
$ws.Range("H77").Value = 37999
$ws.Range("J77").Value = 37999
$ws.Range("L77").Value = 113997
$ws.Range("N77").Value = -122733
$ws.Range("H86").Value = 43425
$ws.Range("I86").Value = 107208.25
$ws.Range("K86").Value = 107208.25
$ws.Range("M86").Value = -106085.25
$ws.Range("H89").Value = 43425
$ws.Range("I89").Value = 107208.25
$ws.Range("K89").Value = 536041.25
$ws.Range("M89").Value = -530425.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1602.4
$ws.Range("I81").Value = 1006
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3018
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -1895
$ws.Range("N81").Value = -8246
$ws.Range("H84").Value = 1602.4
$ws.Range("I84").Value = 1006
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 9054
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -3438
$ws.Range("N84").Value = -29232
$ws.Range("H120").Value = 15824.4
$ws.Range("I120").Value = 10648.8
$ws.Range("K120").Value = 31946.4
$ws.Range("M120").Value = -27108.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 103.36842
$ws.Range("I2").Value = 56.875
$ws.Range("K2").Value = 56.875
$ws.Range("M2").Value = 56.125
$ws.Range("H113").Value = 1470.5834
$ws.Range("I113").Value = 1467.909
$ws.Range("K113").Value = 1467.909
$ws.Range("M113").Value = 702.0909999999999
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H135").Value = 30737920
$ws.Range("J135").Value = 30737920
$ws.Range("L135").Value = 30737920
$ws.Range("N135").Value = -30748060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2772.7334
$ws.Range("I7").Value = 2535.6365
$ws.Range("K7").Value = 2535.6365
$ws.Range("M7").Value = -2423.6365
$ws.Range("H55").Value = 1191.6774
$ws.Range("I55").Value = 1355.9166
$ws.Range("K55").Value = 1355.9166
$ws.Range("M55").Value = -1182.9166
$ws.Range("H68").Value = 5269.65
$ws.Range("I68").Value = 3921.4285
$ws.Range("K68").Value = 3921.4285
$ws.Range("M68").Value = -3172.4285
$ws.Range("H71").Value = 5269.65
$ws.Range("I71").Value = 3921.4285
$ws.Range("K71").Value = 19607.1425
$ws.Range("M71").Value = -15863.1425
$ws.Range("H76").Value = 16798.8
$ws.Range("H79").Value = 16798.8
$ws.Range("H82").Value = 11931.4
$ws.Range("I82").Value = 14290.5625
$ws.Range("K82").Value = 14290.5625
$ws.Range("M82").Value = -13929.5625
$ws.Range("H85").Value = 11931.4
$ws.Range("I85").Value = 14290.5625
$ws.Range("K85").Value = 14290.5625
$ws.Range("M85").Value = -13042.5625
$ws.Range("H122").Value = 6995
$ws.Range("J122").Value = 7995
$ws.Range("L122").Value = 23985
$ws.Range("N122").Value = -28885
$ws.Range("H126").Value = 2772.7334
$ws.Range("I126").Value = 2535.6365
$ws.Range("K126").Value = 7606.9095
$ws.Range("M126").Value = -5136.9095
$ws.Range("H137").Value = 74800
$ws.Range("J137").Value = 74800
$ws.Range("L137").Value = 74800
$ws.Range("N137").Value = -85000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 449.5
$ws.Range("J96").Value = 466
$ws.Range("L96").Value = 466
$ws.Range("N96").Value = -3212
$ws.Range("H100").Value = 1165
$ws.Range("I100").Value = 908
$ws.Range("J100").Value = 2450
$ws.Range("K100").Value = 1816
$ws.Range("L100").Value = 4900
$ws.Range("M100").Value = -1275
$ws.Range("N100").Value = -5982
$ws.Range("H122").Value = 47197.88
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 115294.7
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 345884.1
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -350784.1
$ws.Range("H132").Value = 1680.225
$ws.Range("I132").Value = 1595.1282
$ws.Range("K132").Value = 4785.3846
$ws.Range("M132").Value = -2255.3846
$ws.Range("H136").Value = 69666.94
$ws.Range("I136").Value = 8127.769
$ws.Range("J136").Value = 336336.66
$ws.Range("K136").Value = 24383.307
$ws.Range("L136").Value = 1009009.98
$ws.Range("M136").Value = -21833.307
$ws.Range("N136").Value = -1014109.98
$ws.Range("H137").Value = 58000
$ws.Range("J137").Value = 58000
$ws.Range("L137").Value = 58000
$ws.Range("N137").Value = -68200
